# Update data: 6 November 2021
# Adds the new October 2021 (serial date 44470 = 2021-10-01) unemployment
# data to the "Canada" sheet (one new row) and the "Province" sheet
# (ten new rows, one per province/territory).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada": append a single new row (row 23) for Canada
# ---------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Range("A22:E22").Copy()
$wsCanada.Range("A23:E23").PasteSpecial(-4122)  # xlPasteFormats

$wsCanada.Range("A23").Value = 44470
$wsCanada.Range("B23").Value = "Canada"
$wsCanada.Range("C23").Formula = "=(D23-E23)/E23*100"
$wsCanada.Range("D23").Value = 1365.6
$wsCanada.Range("E23").Value = 1138.4

$wsCanada.Select()
$excel.ActiveWindow.ScrollRow = 16
$wsCanada.Range("A23").Select()

# ---------------------------------------------------------------------
# Sheet "Province": append ten new rows (212-221), one per province,
# for the same reporting date (44470)
# ---------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$wsProvince.Range("A202:E202").Copy()
$wsProvince.Range("A212:E212").PasteSpecial(-4122)  # xlPasteFormats

$wsProvince.Range("A203:E211").Copy()
$wsProvince.Range("A213:E221").PasteSpecial(-4122)  # xlPasteFormats

$provinceRows = @(
    @{ Row = 212; Name = "Newfoundland & Labrador"; D = 35.4;   E = 29.9 },
    @{ Row = 213; Name = "Prince Edward Island";     D = 8;     E = 7.3 },
    @{ Row = 214; Name = "Nova Scotia";              D = 42.1;  E = 41.4 },
    @{ Row = 215; Name = "New Brunswick";            D = 36.4;  E = 31.8 },
    @{ Row = 216; Name = "Quebec";                   D = 255.5; E = 234 },
    @{ Row = 217; Name = "Ontario";                  D = 567.9; E = 422.8 },
    @{ Row = 218; Name = "Manitoba";                 D = 37.1;  E = 36.7 },
    @{ Row = 219; Name = "Saskatchewan";              D = 37.2;  E = 32.1 },
    @{ Row = 220; Name = "Alberta";                  D = 185.9; E = 168.2 },
    @{ Row = 221; Name = "British Columbia";         D = 160;   E = 134.2 }
)

foreach ($r in $provinceRows) {
    $row = $r.Row
    $wsProvince.Range("A$row").Value = 44470
    $wsProvince.Range("B$row").Value = $r.Name
    $wsProvince.Range("C$row").Formula = "=(D$row-E$row)/E$row*100"
    $wsProvince.Range("D$row").Value = $r.D
    $wsProvince.Range("E$row").Value = $r.E
}

$wsProvince.Select()
$excel.ActiveWindow.ScrollRow = 207
$wsProvince.Range("D222").Select()
